$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Project" (F) and "Material" (G) columns for rows 2-7.
# The old placeholder "Material <not specified>" text in column G is
# replaced with a plain number, and column F gets new short labels
# (row 6 becomes a number too).
$ws.Range("F2").Value = "bc"
$ws.Range("G2").Value = 123

$ws.Range("F3").Value = "ad"
$ws.Range("G3").Value = 123

$ws.Range("F4").Value = "qe"
$ws.Range("G4").Value = 123

$ws.Range("F5").Value = "eeee"
$ws.Range("G5").Value = 123

$ws.Range("F6").Value = 111
$ws.Range("G6").Value = 123

$ws.Range("F7").Value = "a"
$ws.Range("G7").Value = 123

# Move the active selection from F10 to G10.
$ws.Range("G10").Select()
